# Adds HPSC data (2020-04-01 .. 2020-04-10) to the "Ireland-manual" sheet
# and restores the selection state, matching the upstream commit
# "Added HPSC data, and updated predictions".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Ireland-JHU
$ws2 = $wb.Worksheets.Item(2)   # Ireland-manual

# New confirmed/death rows for 2020-04-01 .. 2020-04-10, continuing the
# existing layout: B=Country, C=Lat, D=Long, E=date (text), F=cases, G=type
$newRows = @(
    @(142, "2020-04-01", 212, "confirmed"),
    @(143, "2020-04-01", 14,  "death"),
    @(144, "2020-04-02", 402, "confirmed"),
    @(145, "2020-04-02", 13,  "death"),
    @(146, "2020-04-03", 424, "confirmed"),
    @(147, "2020-04-03", 22,  "death"),
    @(148, "2020-04-04", 331, "confirmed"),
    @(149, "2020-04-04", 17,  "death"),
    @(150, "2020-04-05", 390, "confirmed"),
    @(151, "2020-04-05", 21,  "death"),
    @(152, "2020-04-06", 370, "confirmed"),
    @(153, "2020-04-06", 16,  "death"),
    @(154, "2020-04-07", 345, "confirmed"),
    @(155, "2020-04-07", 36,  "death"),
    @(156, "2020-04-08", 365, "confirmed"),
    @(157, "2020-04-08", 25,  "death"),
    @(158, "2020-04-09", 500, "confirmed"),
    @(159, "2020-04-09", 28,  "death"),
    @(160, "2020-04-10", 480, "confirmed"),
    @(161, "2020-04-10", 26,  "death")
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $date   = $r[1]
    $cases  = $r[2]
    $type   = $r[3]

    $ws2.Cells.Item($rowNum, 2).Value = "Ireland"
    $ws2.Cells.Item($rowNum, 3).Value = 53.1424
    $ws2.Cells.Item($rowNum, 4).Value = -7.6921
    # Keep the date column as text (matches the "@" text style already used
    # for column E elsewhere in the sheet) so it isn't coerced to a date serial.
    $ws2.Cells.Item($rowNum, 5).NumberFormat = "@"
    $ws2.Cells.Item($rowNum, 5).Value = $date
    $ws2.Cells.Item($rowNum, 6).Value = $cases
    $ws2.Cells.Item($rowNum, 7).Value = $type
}

# Two trailing blank rows (162/163) - only the date column keeps the text style,
# no values are written (matches the source: <c r="E162" s="1"/>).
$ws2.Cells.Item(162, 5).NumberFormat = "@"
$ws2.Cells.Item(163, 5).NumberFormat = "@"

# "Ireland-manual" becomes/remains the active sheet, with the newly-added
# rows selected and scrolled into view.
$ws2.Activate() | Out-Null
$ws2.Application.Goto($ws2.Range("A162:G163")) | Out-Null
$excel.ActiveWindow.ScrollRow = 133

# "Ireland-JHU" keeps its own selection resting on A2.
$ws1.Activate() | Out-Null
$ws1.Range("A2").Select() | Out-Null

# Restore "Ireland-manual" as the active/selected sheet.
$ws2.Activate() | Out-Null
